$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "27.159.72"
Set-TextValue "E2" "  -1.93%  "
Set-TextValue "D3" "1.556.33"
Set-TextValue "E3" "  -2.17%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "206.22"
Set-TextValue "D6" "0.488"
Set-TextValue "E6" "  -2.43%  "
Set-TextValue "E7" "  +0.04%  "
Set-TextValue "D8" "22.10"
Set-TextValue "E8" "  -0.87%  "
Set-TextValue "E9" "  -2.15%  "
Set-TextValue "D10" "0.0591"
Set-TextValue "E10" "  -0.27%  "
Set-TextValue "E11" "  -0.78%  "
Set-TextValue "D12" "1.779.36"
Set-TextValue "E12" "  -2.05%  "
Set-TextValue "D13" "1.559.79"
Set-TextValue "E13" "  -1.93%  "
Set-TextValue "D14" "3.76"
Set-TextValue "E14" "  -2.53%  "
Set-TextValue "D15" "0.514"
Set-TextValue "E15" "  -3.18%  "
Set-TextValue "D16" "62.81"
Set-TextValue "E16" "  -1.00%  "
Set-TextValue "D17" "27.155.25"
Set-TextValue "E17" "  -1.87%  "
Set-TextValue "D18" "214.18"
Set-TextValue "E18" "  -2.52%  "
Set-TextValue "D19" "0.0₃0684"
Set-TextValue "E19" "  -1.70%  "
Set-TextValue "E20" "  -1.64%  "
Set-TextValue "E21" "  +0.08%  "
Set-TextValue "D22" "4.10"
Set-TextValue "E22" "  -1.05%  "
Set-TextValue "D23" "9.34"
Set-TextValue "E23" "  -3.48%  "
Set-TextValue "E24" "  +0.04%  "
Set-TextValue "D25" "151.93"
Set-TextValue "E25" "  -1.33%  "
Set-TextValue "D26" "6.58"
Set-TextValue "E26" "  -3.30%  "
Set-TextValue "D27" "14.86"
Set-TextValue "E27" "  -1.76%  "
Set-TextValue "E28" "  +0.06%  "
Set-TextValue "E29" "  -1.60%  "
Set-TextValue "E30" "  -1.73%  "
Set-TextValue "E31" "  -1.65%  "
Set-TextValue "E32" "  -1.82%  "
Set-TextValue "D33" "1.376.59"
Set-TextValue "E33" "  -0.01%  "
Set-TextValue "E34" "  -0.61%  "
Set-TextValue "E35" "  -0.40%  "
Set-TextValue "D36" "0.952"
Set-TextValue "E36" "  -2.41%  "
Set-TextValue "E37" "  -1.78%  "
Set-TextValue "E38" "  -1.34%  "
Set-TextValue "E39" "  -2.26%  "
Set-TextValue "D40" "0.514"
Set-TextValue "E40" "  -4.31%  "
Set-TextValue "E41" "  +0.07%  "
Set-TextValue "D42" "0.987"
Set-TextValue "E42" "  +1.88%  "
Set-TextValue "D43" "1.79"
Set-TextValue "E43" "  +3.79%  "
Set-TextValue "E44" "  -0.06%  "
Set-TextValue "E45" "  -1.92%  "
Set-TextValue "E46" "  -0.11%  "
Set-TextValue "D47" "1.691.13"
Set-TextValue "E47" "  -2.04%  "
Set-TextValue "D48" "85.18"
Set-TextValue "E48" "  -2.35%  "
Set-TextValue "D49" "0.0⁷0985"
Set-TextValue "E49" "  -2.46%  "
Set-TextValue "E50" "  -0.60%  "
Set-TextValue "E51" "  +0.17%  "
